# Add a new CRUD record as row 20, following the exact layout of the
# existing data rows (row 19 is the most recent entry and shares every
# text value except the credit number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 19 into row 20 first - this brings along the correct cell
# types (Text for A/B/C/D/E/H, Number for F/G) and keeps every other
# text value ("Victoria Melannye Tibanta Miranda", the RUC/ID numbers,
# "PLANTILLA DE APORTES", the lawyer's name) identical to the source row,
# exactly like the new record in the diff.
$ws.Range("A19:H19").Copy($ws.Range("A20:H20"))

# F20 / G20 are numeric amounts for the new record.
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 23

# A20 needs its own text ("23") instead of the copied "2". Writing the
# string straight into .Value would make Excel re-infer it as a Number,
# losing the Text type the rest of the column uses. Stage it in a scratch
# cell formatted as Text, then copy only the VALUE into A20 so the
# destination's existing (unstyled) format is left untouched.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "23"
$ws.Range("Z1").Copy()
$ws.Range("A20").PasteSpecial("Values")
$ws.Range("Z1").Clear()
